$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the oldest decade (2000年-2009年, rows 2-11). This shifts the
# existing 2010年-2019年 rows up to rows 2-11.
$ws.Rows("2:11").Delete()

# Append the new 2020年 row of data at row 12.
$ws.Range("A12").Value = "2020年"
$ws.Range("B12").Value = 46.2
$ws.Range("C12").Value = 73.7
$ws.Range("D12").Value = 95.3
$ws.Range("E12").Value = 93.09999999999999

# Match the formatting used by the rest of column A (centered, bold,
# bordered "year" label style).
$ws.Range("A11").Copy()
$ws.Range("A12").PasteSpecial(-4122)
$excel.CutCopyMode = 0

Write-Output "done"
